# CdCity.xlsx – replies to SKL inquiries / GenTable updates
#
# Business change: the "刪除" (Delete) marker that was shown in column G
# for the IntRateFloor / IntRateCeiling / IntRateIncr rows (19-21) is
# removed, while the same marker is left in place for the JcicCityCode
# row (22).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")

# Clear the "刪除" marker from G19:G21 (IntRateFloor / IntRateCeiling / IntRateIncr rows)
$ws.Range("G19").Value = $null
$ws.Range("G20").Value = $null
$ws.Range("G21").Value = $null

# Leave G22 (JcicCityCode row) with its "刪除" marker untouched.

# Reflect the author's last on-screen selection before saving.
$ws.Activate()
$ws.Range("I21").Select()
